# Automatic update of files.
# Add a new data row (row 3) to the worksheet, mirroring the existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(3, 1).Value = 111895085
$ws.Cells.Item(3, 2).Value = 90678
$ws.Cells.Item(3, 3).Value = "Ovaliderad"
$ws.Cells.Item(3, 4).Value = "LC"
$ws.Cells.Item(3, 5).Value = 4366
$ws.Cells.Item(3, 6).Value = "Skarp dropptaggsvamp"
$ws.Cells.Item(3, 7).Value = "Hydnellum peckii"
$ws.Cells.Item(3, 8).Value = "Banker"
$ws.Cells.Item(3, 9).Value = ""
$ws.Cells.Item(3, 11).Value = ""
$ws.Cells.Item(3, 16).Value = "Kratte masugn (Kratte masugn), Gstr"
$ws.Cells.Item(3, 17).Value = 576346.0152053731
$ws.Cells.Item(3, 18).Value = 6702381.515453912
$ws.Cells.Item(3, 19).Value = 1
$ws.Cells.Item(3, 20).Value = "Gävleborg"
$ws.Cells.Item(3, 21).Value = "Hofors"
$ws.Cells.Item(3, 22).Value = "Gästrikland"
$ws.Cells.Item(3, 23).Value = "Torsåker"
$ws.Cells.Item(3, 25).NumberFormat = "@"
$ws.Cells.Item(3, 25).Value = "2023-09-04"
$ws.Cells.Item(3, 25).Style = "Normal"
$ws.Cells.Item(3, 26).Value = "17:38"
$ws.Cells.Item(3, 27).NumberFormat = "@"
$ws.Cells.Item(3, 27).Value = "2023-09-04"
$ws.Cells.Item(3, 27).Style = "Normal"
$ws.Cells.Item(3, 28).Value = "17:38"
$ws.Cells.Item(3, 30).Value = $false
$ws.Cells.Item(3, 31).Value = $false
$ws.Cells.Item(3, 33).Value = $false
$ws.Cells.Item(3, 46).Value = ""
$ws.Cells.Item(3, 49).Value = "fanny westling"
$ws.Cells.Item(3, 50).Value = "fanny westling"
$ws.Cells.Item(3, 51).Value = ""
